# Refresh cryptos list: prices (Price) and 1h volume deltas (Volume(1h))
# Mirrors a scheduled GitHub Actions data refresh of cryptos.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.148.28'
$ws.Range('E2').Value = '  -1.85%  '

$ws.Range('D3').Value = '1.558.58'
$ws.Range('E3').Value = '  -1.84%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').Value = '''206.49'
$ws.Range('E5').Value = '  -0.31%  '

$ws.Range('E6').Value = '  -1.19%  '

$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').Value = '''22.26'
$ws.Range('E8').Value = '  +0.16%  '

$ws.Range('E9').Value = '  -1.88%  '

$ws.Range('E10').Value = '  +0.18%  '

$ws.Range('D11').Value = '''0.0862'
$ws.Range('E11').Value = '  -0.62%  '

$ws.Range('D12').Value = '1.782.36'
$ws.Range('E12').Value = '  -1.72%  '

$ws.Range('D13').Value = '1.558.03'
$ws.Range('E13').Value = '  -1.96%  '

$ws.Range('D14').Value = '''3.76'
$ws.Range('E14').Value = '  -2.32%  '

$ws.Range('D15').Value = '''0.514'
$ws.Range('E15').Value = '  -3.00%  '

$ws.Range('D16').Value = '''62.78'
$ws.Range('E16').Value = '  -1.05%  '

$ws.Range('D17').Value = '27.153.17'
$ws.Range('E17').Value = '  -1.81%  '

$ws.Range('D18').Value = '''213.62'
$ws.Range('E18').Value = '  -2.65%  '

$ws.Range('D19').Value = '0.0₃0686'
$ws.Range('E19').Value = '  -1.30%  '

$ws.Range('E20').Value = '  -1.23%  '

$ws.Range('D22').Value = '''4.10'
$ws.Range('E22').Value = '  -0.75%  '

$ws.Range('E23').Value = '  -3.12%  '

$ws.Range('E24').Value = '  -0.03%  '

$ws.Range('D25').Value = '''151.91'
$ws.Range('E25').Value = '  -1.11%  '

$ws.Range('D26').Value = '''6.59'
$ws.Range('E26').Value = '  -3.28%  '

$ws.Range('D27').Value = '''14.87'
$ws.Range('E27').Value = '  -1.66%  '

$ws.Range('E28').Value = '  +0.02%  '

$ws.Range('E29').Value = '  -1.20%  '

$ws.Range('D30').Value = '''1.14'
$ws.Range('E30').Value = '  -0.58%  '

$ws.Range('E31').Value = '  -1.22%  '

$ws.Range('D32').Value = '''3.16'
$ws.Range('E32').Value = '  -1.96%  '

$ws.Range('D33').Value = '1.382.42'
$ws.Range('E33').Value = '  +0.96%  '

$ws.Range('E34').Value = '  +0.29%  '

$ws.Range('E35').Value = '  +0.16%  '

$ws.Range('D36').Value = '''0.947'
$ws.Range('E36').Value = '  -2.81%  '

$ws.Range('D38').Value = '''0.0165'
$ws.Range('E38').Value = '  -1.32%  '

$ws.Range('D39').Value = '''0.813'
$ws.Range('E39').Value = '  -1.44%  '

$ws.Range('D40').Value = '''0.516'
$ws.Range('E40').Value = '  -3.57%  '

$ws.Range('D42').Value = '''0.986'
$ws.Range('E42').Value = '  +1.48%  '

$ws.Range('E43').Value = '  +3.52%  '

$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '''63.27'
$ws.Range('E44').Value = '  -1.38%  '

$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').Value = '''2.16'
$ws.Range('E45').Value = '  -0.20%  '

$ws.Range('D46').Value = '''5.21'
$ws.Range('E46').Value = '  +0.64%  '

$ws.Range('D47').Value = '1.694.11'
$ws.Range('E47').Value = '  -1.73%  '

$ws.Range('D48').Value = '''85.44'
$ws.Range('E48').Value = '  -2.23%  '

$ws.Range('D49').Value = '0.0₇0994'
$ws.Range('E49').Value = '  -0.91%  '

$ws.Range('E50').Value = '  -0.34%  '

$ws.Range('E51').Value = '  +0.13%  '
